$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 2) so that shared strings for the IMF/IMF(20%) columns are correct
$ws.Range("D2").Value = "IMF (20%) - Sales"
$ws.Range("E2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("F2").Value = "IMF - Sales"
$ws.Range("G2").Value = "IMF - Sales + Emp"
$ws.Range("L2").Value = "IMF (20%) - Sales"
$ws.Range("M2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("N2").Value = "IMF - Sales"
$ws.Range("O2").Value = "IMF - Sales + Emp"
$ws.Range("T2").Value = "IMF (20%) - Sales"
$ws.Range("U2").Value = "IMF (20%) - Sales + Emp"
$ws.Range("V2").Value = "IMF - Sales"
$ws.Range("W2").Value = "IMF - Sales + Emp"

# Update data cells rows 4-13 for columns D,E,F,G / L,M,N,O / T,U,V,W

# Row 4
$ws.Range("D4").Value = 0.9242685904175726
$ws.Range("E4").Value = 0.9396081541203916
$ws.Range("F4").Value = 4.621342952087861
$ws.Range("G4").Value = 4.698040770601957
$ws.Range("L4").Value = 0.2939782475134149
$ws.Range("M4").Value = 0.2935480512885106
$ws.Range("N4").Value = 0.2939782475134149
$ws.Range("O4").Value = 0.2935480512885106
$ws.Range("T4").Value = 60206736130
$ws.Range("U4").Value = 60380413875
$ws.Range("V4").Value = 60206736130
$ws.Range("W4").Value = 60380413875

# Row 5
$ws.Range("D5").Value = 0.7383447073709195
$ws.Range("E5").Value = 0.8065494146227172
$ws.Range("F5").Value = 3.691723536854596
$ws.Range("G5").Value = 4.032747073113593
$ws.Range("L5").Value = 0.2187022724763601
$ws.Range("M5").Value = 0.2201697804353018
$ws.Range("N5").Value = 0.2187022724763601
$ws.Range("O5").Value = 0.2201697804353018
$ws.Range("T5").Value = 976531986457
$ws.Range("U5").Value = 988562844368
$ws.Range("V5").Value = 976531986457
$ws.Range("W5").Value = 988562844368

# Row 6
$ws.Range("D6").Value = 0.1800287566322612
$ws.Range("E6").Value = 0.5127070547825061
$ws.Range("F6").Value = 0.9001437831613054
$ws.Range("G6").Value = 2.56353527391253
$ws.Range("L6").Value = 0.3972232104675282
$ws.Range("M6").Value = 0.3787799051437756
$ws.Range("N6").Value = 0.3972232104675282
$ws.Range("O6").Value = 0.3787799051437756
$ws.Range("T6").Value = 25762595315
$ws.Range("U6").Value = 37793453226
$ws.Range("V6").Value = 25762595315
$ws.Range("W6").Value = 37793453226

# Row 7
$ws.Range("D7").Value = 0.3035777241591624
$ws.Range("E7").Value = 0.7698259819696905
$ws.Range("F7").Value = 1.517888620795812
$ws.Range("G7").Value = 3.849129909848453
$ws.Range("L7").Value = 0.5468907065540652
$ws.Range("M7").Value = 0.4790830303318274
$ws.Range("N7").Value = 0.5468907065540652
$ws.Range("O7").Value = 0.4790830303318274
$ws.Range("T7").Value = 32205923465
$ws.Range("U7").Value = 46231904709
$ws.Range("V7").Value = 32205923465
$ws.Range("W7").Value = 46231904709

# Row 8
$ws.Range("D8").Value = 1.260193722463467
$ws.Range("E8").Value = 0.9921647696433968
$ws.Range("F8").Value = 6.300968612317333
$ws.Range("G8").Value = 4.960823848216997
$ws.Range("L8").Value = 0.2083789632708606
$ws.Range("M8").Value = 0.2083789632708606
$ws.Range("N8").Value = 0.2083789632708606
$ws.Range("O8").Value = 0.2083789632708606
$ws.Range("T8").Value = 933471841988
$ws.Range("U8").Value = 933471841988
$ws.Range("V8").Value = 933471841988
$ws.Range("W8").Value = 933471841988

# Row 9
$ws.Range("D9").Value = 1.238354713022687
$ws.Range("E9").Value = 2.005485635699418
$ws.Range("F9").Value = 6.191773565113434
$ws.Range("G9").Value = 10.02742817849709
$ws.Range("L9").Value = 1.348214135245821
$ws.Range("M9").Value = 0.5441963828111374
$ws.Range("N9").Value = 1.348214135245821
$ws.Range("O9").Value = 0.5441963828111374
$ws.Range("T9").Value = 6548180327
$ws.Range("U9").Value = 26602310262
$ws.Range("V9").Value = 6548180327
$ws.Range("W9").Value = 26602310262

# Row 10
$ws.Range("D10").Value = 0.9868244542315268
$ws.Range("E10").Value = 1.073704801281935
$ws.Range("F10").Value = 4.934122271157631
$ws.Range("G10").Value = 5.368524006409673
$ws.Range("L10").Value = 0.2733403708402601
$ws.Range("M10").Value = 0.2733403708402601
$ws.Range("N10").Value = 0.2733403708402601
$ws.Range("O10").Value = 0.2733403708402601
$ws.Range("T10").Value = 65619795685
$ws.Range("U10").Value = 65619795685
$ws.Range("V10").Value = 65619795685
$ws.Range("W10").Value = 65619795685

# Row 11
$ws.Range("D11").Value = 1.18860944939988
$ws.Range("E11").Value = 1.053520392819479
$ws.Range("F11").Value = 5.943047246999395
$ws.Range("G11").Value = 5.267601964097405
$ws.Range("L11").Value = 0.2202869195466604
$ws.Range("M11").Value = 0.2202731665795292
$ws.Range("N11").Value = 0.2202869195466604
$ws.Range("O11").Value = 0.2202731665795292
$ws.Range("T11").Value = 957938292874
$ws.Range("U11").Value = 958111970619
$ws.Range("V11").Value = 957938292874
$ws.Range("W11").Value = 958111970619

# Row 12
$ws.Range("D12").Value = 1.014270909536257
$ws.Range("E12").Value = 1.128821247125711
$ws.Range("F12").Value = 5.071354547681283
$ws.Range("G12").Value = 5.644106235628553
$ws.Range("L12").Value = 0.3823216259415932
$ws.Range("M12").Value = 0.3818638512960749
$ws.Range("N12").Value = 0.3823216259415932
$ws.Range("O12").Value = 0.3818638512960749
$ws.Range("T12").Value = 90086246571
$ws.Range("U12").Value = 90259924316
$ws.Range("V12").Value = 90086246571
$ws.Range("W12").Value = 90259924316

# Row 13
$ws.Range("D13").Value = 1.486349685879603
$ws.Range("E13").Value = 0.9247074269381814
$ws.Range("F13").Value = 7.431748429398017
$ws.Range("G13").Value = 4.623537134690935
$ws.Range("L13").Value = 0.2034671187804626
$ws.Range("M13").Value = 0.2034671187804626
$ws.Range("N13").Value = 0.2034671187804626
$ws.Range("O13").Value = 0.2034671187804626
$ws.Range("T13").Value = 867852046303
$ws.Range("U13").Value = 867852046303
$ws.Range("V13").Value = 867852046303
$ws.Range("W13").Value = 867852046303
